$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column W is column index 23 (A=1). The values in W represent a
# "children" category (0-5) that was being mis-calculated; the fix maps
# each non-zero category to its corrected value: new = old*10 + 15
# (0 stays 0, 1->25, 2->35, 3->45, 4->55, 5->65).
$lastRow = $ws.Cells.SpecialCells(11).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 23)
    $old = $cell.Value2
    if ($old -ne 0) {
        $cell.Value = $old * 10 + 15
    }
}
